# The author retyped part of "rackets" -> "rоckets" (the "a" was
# replaced by a Cyrillic "о", U+043E) while the cursor/insertion point
# sat right at the "_GoBack" bookmark that separates "them developing"
# from " rackets and why not to invent". As a result the final OOXML
# shows:
#   ...them developing</w:r>
#   <w:r> rо</w:r>                 (new run, inserted before the bookmark)
#   <w:bookmarkStart.../><w:bookmarkEnd.../>
#   <w:r>ckets and why not to invent</w:r>   (old run, lost its " ra" prefix)
#
# Net visible effect: " rackets" -> " rоckets".

$d = $word.ActiveDocument

$bm = $d.Bookmarks("_GoBack")
$bmStart = $bm.Start

# 1. Insert " r" + Cyrillic "о" immediately before the bookmark.
$insPoint = $d.Range($bmStart, $bmStart)
$insPoint.InsertBefore(" r" + [char]0x043E)

# 2. Force the newly inserted text to live in its own run instead of
#    being silently re-merged into the preceding "them developing" run:
#    toggle a formatting property on and back off.
$newRange = $d.Range($bmStart, $bmStart + 3)
$newRange.Font.Bold = 1
$newRange.Font.Bold = 0

# 3. Remove the now-duplicated " ra" that used to start " rackets..."
#    (it immediately follows the bookmark, which has not moved).
$bm = $d.Bookmarks("_GoBack")
$oldRunStart = $bm.End
$delRange = $d.Range($oldRunStart, $oldRunStart + 3)
$delRange.Text = ""

# 4. The deletion above causes this COM host to silently coalesce the
#    "ckets and why not to invent" run forward with the still-untouched,
#    identically-formatted runs that follow it in the same paragraph
#    (" flying cars.", " For now I want to grow as", " Front-end
#    developer."), which the source diff does NOT touch. Re-establish
#    those original run boundaries by nudging formatting on/off again.
$r1 = $d.Content
$r1.Find.Execute(" flying cars.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$flyingRange = $d.Range($r1.Start, $r1.End)
$flyingRange.Font.Bold = 1
$flyingRange.Font.Bold = 0

$r2 = $d.Content
$r2.Find.Execute(" For now I want to grow as", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$forNowRange = $d.Range($r2.Start, $r2.End)
$forNowRange.Font.Bold = 1
$forNowRange.Font.Bold = 0

$r3 = $d.Content
$r3.Find.Execute(" Front-end developer.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$frontEndRange = $d.Range($r3.Start, $r3.End)
$frontEndRange.Font.Bold = 1
$frontEndRange.Font.Bold = 0
